# "Report algo global pheromone maj"
#
# The "ΔΣT comparatif" summary block (rows 1:35, columns A:E) on Sheet1
# drops its two "Δv" / "Δw" columns (D and E), keeping only the
# B ("ΔΣT") and C ("Δtemps") comparison columns. A new summary row is
# added underneath (row 36, "%Δ") expressing the MOY swing as a
# percentage of the MIN/MAX extremes. The remembered selection also
# moves from F37 to F27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Δv" / "Δw" columns from the comparison block:
#  - D2:E2    header labels
#  - D3:E32   per-row ΔΣT/Δtemps/Δv/Δw comparisons
#  - D33:E35  MIN / MAX / MOY summary formulas
$ws.Range("D2:E2").ClearContents()
$ws.Range("D3:E32").ClearContents()
$ws.Range("D33:E35").ClearContents()

# New row: "%Δ" of the MOY relative to the MAX (col B) / MIN (col C) extremes.
$ws.Range("A36").Value = "%Δ"
$ws.Range("B36").Formula = "=(B35/B34)*100"
$ws.Range("C36").Formula = "=C35/C33*100"

# Selection moved from F37 to F27.
[void]$ws.Range("F27").Select()
